$d = $word.ActiveDocument

# 1. Update the date
$d.Content.Find.Execute("August 18, 2020", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "August 20, 2020", 2)

# 2. Update the address
$d.Content.Find.Execute("6/F Filinvest Bldg., No. 79 EDSA, Highway Hills, Mandaluyong City", `
                         $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Brgy. Kilada, Matalam, Cotabato", 2)

# 3. Update the production estimate / fee paragraph
$d.Content.Find.Execute( `
    "Please be informed that based on your submitted production estimate of 1,000.00 Metric Tons or 20,000.00 Lkg., your Milling License Fee for Crop Year 2020 - 2021 is ONE THOUSAND  (PHP 1,000.00) PESOS.  However, you have an underpayment in your Milling License Fee for CY 2020 - 2021 in the amount of ONE THOUSAND  PESOS (PHP 1,000.00).", `
    $false, $false, $false, $false, $false, `
    $true, 1, $false, `
    "Please be informed that based on your submitted production estimate of 0.00 Metric Tons or 0.00 Lkg., your Milling License Fee for Crop Year 2020 - 2021 is  (PHP 0.00) PESOS.  However, you have an excess payment in your Milling License Fee for CY 2020 - 2021 in the amount of  PESOS (PHP 0.00).", `
    2)

# 4. Update the bolded amount due
$d.Content.Find.Execute("TWO THOUSAND  PESOS (PHP 2,000.00)", $false, $false, $false, $false, $false, `
                         $true, 1, $false, " PESOS (PHP 0.00)", 2)
